# Update transition-probability matrix on Sheet1 with refreshed values
# from games pulled March 7 (Kansas City_B team-specific matrix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1769436997319035
$ws.Range("C2").Value = 0.5844504021447721
$ws.Range("J2").Value = 0.01876675603217158
$ws.Range("P2").Value = 0.1313672922252011
$ws.Range("S2").Value = 0.08847184986595175
$ws.Range("B3").Value = 0.009009009009009009
$ws.Range("C3").Value = 0.03153153153153153
$ws.Range("J3").Value = 0.02252252252252252
$ws.Range("P3").Value = 0.6891891891891891
$ws.Range("S3").Value = 0.2477477477477477
$ws.Range("J4").Value = 0.05357142857142857
$ws.Range("P4").Value = 0.6428571428571429
$ws.Range("S4").Value = 0.3035714285714285
$ws.Range("B6").Value = 0.1260162601626016
$ws.Range("D6").Value = 0.02439024390243903
$ws.Range("F6").Value = 0.07723577235772358
$ws.Range("J6").Value = 0.2520325203252032
$ws.Range("O6").Value = 0.03252032520325204
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.02845528455284553
$ws.Range("S6").Value = 0.2926829268292683
$ws.Range("B7").Value = 0.1207729468599034
$ws.Range("D7").Value = 0.01449275362318841
$ws.Range("F7").Value = 0.05797101449275362
$ws.Range("J7").Value = 0.1400966183574879
$ws.Range("O7").Value = 0.03864734299516908
$ws.Range("Q7").Value = 0.1014492753623188
$ws.Range("R7").Value = 0.06280193236714976
$ws.Range("S7").Value = 0.463768115942029
$ws.Range("B8").Value = 0.1230769230769231
$ws.Range("D8").Value = 0.02637362637362637
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.05934065934065934
$ws.Range("J8").Value = 0.1054945054945055
$ws.Range("O8").Value = 0.01758241758241758
$ws.Range("Q8").Value = 0.167032967032967
$ws.Range("R8").Value = 0.09010989010989011
$ws.Range("S8").Value = 0.4087912087912088
$ws.Range("B9").Value = 0.1101321585903084
$ws.Range("D9").Value = 0.03083700440528634
$ws.Range("F9").Value = 0.07488986784140969
$ws.Range("J9").Value = 0.1101321585903084
$ws.Range("O9").Value = 0.00881057268722467
$ws.Range("Q9").Value = 0.1806167400881057
$ws.Range("R9").Value = 0.0881057268722467
$ws.Range("S9").Value = 0.3964757709251101
$ws.Range("B10").Value = 0.1268715524034673
$ws.Range("D10").Value = 0.02364066193853428
$ws.Range("E10").Value = 0.003152088258471237
$ws.Range("F10").Value = 0.07013396375098503
$ws.Range("J10").Value = 0.1245074862096139
$ws.Range("O10").Value = 0.02048857368006304
$ws.Range("Q10").Value = 0.1954294720252167
$ws.Range("R10").Value = 0.07171000788022065
$ws.Range("S10").Value = 0.3640661938534279
$ws.Range("G11").Value = 0.1779935275080906
$ws.Range("J11").Value = 0.07119741100323625
$ws.Range("K11").Value = 0.2233009708737864
$ws.Range("L11").Value = 0.511326860841424
$ws.Range("S11").Value = 0.01618122977346278
$ws.Range("G12").Value = 0.7891566265060241
$ws.Range("J12").Value = 0.1385542168674699
$ws.Range("L12").Value = 0.03614457831325301
$ws.Range("S12").Value = 0.03614457831325301
$ws.Range("G13").Value = 0.6486486486486487
$ws.Range("J13").Value = 0.2162162162162162
$ws.Range("S13").Value = 0.1351351351351351
$ws.Range("F15").Value = 0.02508960573476703
$ws.Range("H15").Value = 0.1397849462365591
$ws.Range("I15").Value = 0.09677419354838709
$ws.Range("J15").Value = 0.3405017921146953
$ws.Range("K15").Value = 0.05734767025089606
$ws.Range("M15").Value = 0.007168458781362007
$ws.Range("O15").Value = 0.08960573476702509
$ws.Range("S15").Value = 0.2437275985663082
$ws.Range("F16").Value = 0.02155172413793104
$ws.Range("H16").Value = 0.1724137931034483
$ws.Range("I16").Value = 0.08189655172413793
$ws.Range("J16").Value = 0.3448275862068966
$ws.Range("K16").Value = 0.09913793103448276
$ws.Range("M16").Value = 0.01724137931034483
$ws.Range("N16").Value = 0.004310344827586207
$ws.Range("O16").Value = 0.1120689655172414
$ws.Range("S16").Value = 0.146551724137931
$ws.Range("F17").Value = 0.02132701421800948
$ws.Range("H17").Value = 0.1682464454976303
$ws.Range("I17").Value = 0.09004739336492891
$ws.Range("J17").Value = 0.4194312796208531
$ws.Range("K17").Value = 0.0947867298578199
$ws.Range("M17").Value = 0.01184834123222749
$ws.Range("O17").Value = 0.08767772511848342
$ws.Range("S17").Value = 0.1066350710900474
$ws.Range("F18").Value = 0.03448275862068965
$ws.Range("H18").Value = 0.1609195402298851
$ws.Range("I18").Value = 0.1264367816091954
$ws.Range("J18").Value = 0.3908045977011494
$ws.Range("K18").Value = 0.09770114942528736
$ws.Range("M18").Value = 0.01149425287356322
$ws.Range("O18").Value = 0.05747126436781609
$ws.Range("S18").Value = 0.1206896551724138
$ws.Range("F19").Value = 0.02225519287833828
$ws.Range("H19").Value = 0.2010385756676558
$ws.Range("I19").Value = 0.09050445103857567
$ws.Range("J19").Value = 0.3605341246290801
$ws.Range("K19").Value = 0.1060830860534125
$ws.Range("M19").Value = 0.01780415430267062
$ws.Range("O19").Value = 0.0712166172106825
$ws.Range("S19").Value = 0.1305637982195846
